$d = $word.ActiveDocument

# Update the date line at top of document
$d.Content.Find.Execute("2025-03-05 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-06 Thursday", 2)

# Update the multiplication problems in the table.
# Using direct cell access (row, col) avoids ambiguity from find/replace
# when a new value coincides with an old value found elsewhere in the table.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "99×44="
$t.Cell(1, 2).Range.Text = "78×35="
$t.Cell(1, 3).Range.Text = "34×78="
$t.Cell(1, 4).Range.Text = "99×69="
$t.Cell(1, 5).Range.Text = "67×55="

$t.Cell(5, 1).Range.Text = "42×25="
$t.Cell(5, 2).Range.Text = "13×59="
$t.Cell(5, 3).Range.Text = "57×58="
$t.Cell(5, 4).Range.Text = "33×20="
$t.Cell(5, 5).Range.Text = "87×83="

$t.Cell(10, 1).Range.Text = "72×92="
$t.Cell(10, 2).Range.Text = "24×75="
$t.Cell(10, 3).Range.Text = "18×56="
$t.Cell(10, 4).Range.Text = "46×49="
$t.Cell(10, 5).Range.Text = "66×11="

$t.Cell(15, 1).Range.Text = "15×48="
$t.Cell(15, 2).Range.Text = "60×48="
$t.Cell(15, 3).Range.Text = "98×45="
$t.Cell(15, 4).Range.Text = "47×86="
$t.Cell(15, 5).Range.Text = "57×66="

$t.Cell(20, 1).Range.Text = "47×76="
$t.Cell(20, 2).Range.Text = "64×11="
$t.Cell(20, 3).Range.Text = "83×36="
$t.Cell(20, 4).Range.Text = "48×47="
$t.Cell(20, 5).Range.Text = "67×16="
